$d = $word.ActiveDocument

# 1. Investment accounts bullet
$d.Content.Find.Execute(
    "Worked with financial accounts to produce monetary growth through Mutual Funds, Swing Trading, and Value Investing. Charted and analyzed with code and excel to see top performing stocks and strategies.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Managed accounts to produce monetary growth through Mutual Funds, Swing Trading, and Value Investing. Tracked and analyzed monthly data with Excel to highlight highest performing stocks and strategies.",
    2
)

# 2. Central Park instructor bullet
$d.Content.Find.Execute(
    "Led daily outdoor and art based classes for elementary aged students in Central Park and afterwards helped them understand academic material they were learning in school.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Instructed daily outdoor and art based classes of ten or more elementary aged students in Central Park. Furthermore, I helped them in understanding their in-school academic material.",
    2
)

# 3. Peer helper bullet
$d.Content.Find.Execute(
    "Helped peers with multiple computer languages (Python, C++, Unix) to succeed in the class.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Helped peers in labs of 15-30 students with multiple computer languages (Python, C++, Unix) to succeed in the class.",
    2
)

# 4. Languages paragraph: text + formatting
$d.Content.Find.Execute(
    "English, French, Russian. Learning Spanish and Bangla",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "English, French, Russian. Learning Spanish and Bangla.",
    2
)

foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "English, French, Russian*") {
        $p.Range.Font.Italic = $true
        $p.LineSpacingRule = 5   # wdLineSpaceExactly? use explicit below instead
        $p.Format.LineSpacing = 18
        break
    }
}
